# Updates the cryptocurrency list data (coin names/links/prices/1h volume %)
# to reflect the "Updated cryptos list ... with GitHub Actions" refresh.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value is a plain decimal number (e.g. "226.82") must be
# forced to stay as literal text, matching the source inlineStr cells,
# otherwise Excel's automatic type inference would convert them to numbers
# (losing formatting like trailing zeros, e.g. "4.20" -> 4.2).
function Set-TextValue($addr, $value) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

# Row 2
$ws.Range("D2").Value = '34.056.39'
$ws.Range("E2").Value = '  -0.16%  '
# Row 3
$ws.Range("D3").Value = '1.788.84'
$ws.Range("E3").Value = '  -0.05%  '
# Row 4
$ws.Range("E4").Value = '  +0.04%  '
# Row 5
Set-TextValue "D5" '226.82'
$ws.Range("E5").Value = '  +1.82%  '
# Row 6
$ws.Range("E6").Value = '  -1.35%  '
# Row 7
$ws.Range("E7").Value = '  -0.04%  '
# Row 8
Set-TextValue "D8" '32.29'
$ws.Range("E8").Value = '  -0.20%  '
# Row 9
$ws.Range("E9").Value = '  +3.95%  '
# Row 10
Set-TextValue "D10" '0.0684'
$ws.Range("E10").Value = '  -4.47%  '
# Row 11
Set-TextValue "D11" '0.0940'
$ws.Range("E11").Value = '  +1.09%  '
# Row 12
$ws.Range("D12").Value = '2.046.58'
$ws.Range("E12").Value = '  -0.03%  '
# Row 13
Set-TextValue "D13" '11.44'
$ws.Range("E13").Value = '  +4.45%  '
# Row 14
$ws.Range("D14").Value = '1.786.59'
$ws.Range("E14").Value = '  -0.30%  '
# Row 15
Set-TextValue "D15" '0.622'
$ws.Range("E15").Value = '  -0.71%  '
# Row 16
$ws.Range("D16").Value = '34.029.88'
$ws.Range("E16").Value = '  -0.16%  '
# Row 17
Set-TextValue "D17" '4.20'
$ws.Range("E17").Value = '  +0.53%  '
# Row 18
Set-TextValue "D18" '67.97'
$ws.Range("E18").Value = '  -0.15%  '
# Row 19
Set-TextValue "D19" '242.67'
$ws.Range("E19").Value = '  -0.67%  '
# Row 20
$ws.Range("E20").Value = '  -1.46%  '
# Row 21
$ws.Range("E21").Value = '  +0.04%  '
# Row 22
$ws.Range("E22").Value = '  -0.17%  '
# Row 23
Set-TextValue "D23" '4.09'
$ws.Range("E23").Value = '  +0.00%  '
# Row 24
Set-TextValue "D24" '2.06'
$ws.Range("E24").Value = '  -2.76%  '
# Row 25
Set-TextValue "D25" '162.12'
$ws.Range("E25").Value = '  +2.09%  '
# Row 26
$ws.Range("E26").Value = '  +1.08%  '
# Row 27
$ws.Range("E27").Value = '  -0.96%  '
# Row 28
$ws.Range("E28").Value = '  +0.51%  '
# Row 29
$ws.Range("E29").Value = '  +0.12%  '
# Row 30
$ws.Range("E30").Value = '  +2.71%  '
# Row 31
$ws.Range("E31").Value = '  -0.65%  '
# Row 32
$ws.Range("E32").Value = '  -0.83%  '
# Row 33
$ws.Range("E33").Value = '  +3.22%  '
# Row 34
Set-TextValue "D34" '1.84'
$ws.Range("E34").Value = '  +1.81%  '
# Row 35
$ws.Range("D35").Value = '1.398.64'
$ws.Range("E35").Value = '  +0.29%  '
# Row 36
Set-TextValue "D36" '0.652'
$ws.Range("E36").Value = '  +0.37%  '
# Row 37
$ws.Range("E37").Value = '  -0.72%  '
# Row 38
$ws.Range("E38").Value = '  +8.95%  '
# Row 39
$ws.Range("E39").Value = '  +1.50%  '
# Row 40
Set-TextValue "D40" '80.10'
$ws.Range("E40").Value = '  +0.46%  '
# Row 41
$ws.Range("E41").Value = '  +0.04%  '
# Row 42
Set-TextValue "D42" '0.919'
$ws.Range("E42").Value = '  -0.10%  '
# Row 43
Set-TextValue "D43" '13.72'
$ws.Range("E43").Value = '  +14.41%  '
# Row 44
$ws.Range("E44").Value = '  -0.51%  '
# Row 45
$ws.Range("B45").Value = 'FraxShare'
$ws.Range("C45").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
Set-TextValue "D45" '6.10'
$ws.Range("E45").Value = '  +2.24%  '
# Row 46
$ws.Range("B46").Value = 'BabyDogeCoin'
$ws.Range("C46").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D46").Value = '0.0₆0139'
$ws.Range("E46").Value = '  +8.16%  '
# Row 47
$ws.Range("B47").Value = 'WEMIXToken'
$ws.Range("C47").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
Set-TextValue "D47" '1.08'
$ws.Range("E47").Value = '  +2.67%  '
# Row 48
$ws.Range("B48").Value = 'Kaspa'
$ws.Range("C48").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
Set-TextValue "D48" '0.0506'
$ws.Range("E48").Value = '  +1.79%  '
# Row 49
Set-TextValue "D49" '107.68'
$ws.Range("E49").Value = '  +0.12%  '
# Row 50
$ws.Range("D50").Value = '1.947.36'
$ws.Range("E50").Value = '  -0.14%  '
# Row 51
$ws.Range("E51").Value = '  +0.04%  '
